$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 19 (push it down to row 21).
$ws.Rows("19:20").Insert()

# New row 19: numeric weather entry (Cloudy, 17).
$ws.Range("A19").Value = 26
$ws.Range("B19").Value = "Cloudy"
$ws.Range("C19").Value = "'01/08/2025"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = 17

# New row 20: numeric weather entry (Rain, 22).
$ws.Range("A20").Value = 26
$ws.Range("B20").Value = "Rain"
$ws.Range("C20").Value = "'01/08/2025"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = 22

# Row 21 is the original (text-typed) row, pushed down; only its D value changes (17 -> 22),
# keeping the original text type.
$ws.Range("D21").Value = "'22"
$ws.Range("D21").Style = "Normal"
